$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 16; this shifts the existing rows 16..84
# down to 17..85, matching the target diff (which appends a new weekly
# record and pushes the rest of the daily log down by one row).
$ws.Rows("16:16").Insert()

# Populate the newly inserted row 16 with the new record's data.
$ws.Range("A16").Value = 11
$ws.Range("B16").Value = "Vega Monumental Concepción"
$ws.Range("C16").Value = "Bíobío"
$ws.Range("D16").Value = 44623
$ws.Range("E16").Value = 8
$ws.Range("F16").Value = 100112024
$ws.Range("G16").Value = "Choclo"
$ws.Range("H16").Value = "Choclero"
$ws.Range("I16").Value = "Primera"
$ws.Range("J16").Value = 3500
$ws.Range("K16").Value = 200
$ws.Range("L16").Value = 200
$ws.Range("M16").Value = 200
$ws.Range("N16").Value = "$/unidad"
$ws.Range("O16").Value = "Región Metropolitana"
$ws.Range("P16").Value = 200
$ws.Range("Q16").Value = 1
$ws.Range("R16").Value = "Hortaliza"
